$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(61919, "Alícia Rocha", "Engenharia", "Problemas pessoais", 7, 45081, 8859.82),
    @(85614, "Maria Clara Gonçalves", "Recursos Humanos", "Doença", 1, 45090, 11113.23),
    @(58994, "Diogo das Neves", "Marketing", "Consulta médica", 7, 45081, 12161.9),
    @(6038, "Bruna Moura", "Recursos Humanos", "Consulta médica", 2, 45105, 8812.18),
    @(94586, "Maria Eduarda Gonçalves", "P&D", "Consulta médica", 4, 45083, 5416.99),
    @(26379, "Raquel Nogueira", "Recursos Humanos", "Problemas pessoais", 6, 45105, 5454.48),
    @(58593, "Luigi Costela", "Operações", "Consulta médica", 5, 45081, 8756.18),
    @(32985, "Maysa Cavalcanti", "Recursos Humanos", "Outros", 8, 45093, 3663),
    @(15615, "Davi Pires", "Financeiro", "Consulta médica", 7, 45100, 10970.53),
    @(12645, "Vitória Costela", "Recursos Humanos", "Consulta médica", 3, 45086, 5219.13)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
